$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("C2").Value = "chromegrid"
$ws.Range("C3").Value = "android"
$ws.Range("C4").Value = "chromegrid"
$ws.Range("C3").Select()
